$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2: update reporting date serial value
$ws.Range("B2").Value = 45758.0

# Rows 6-17 and 19-47: convert shared formula column E into static rounded
# values, and round column F values to the nearest whole number.
$rowData = @(
    @(6, 0.73, 68765.0),
    @(7, 0.81, 75023.0),
    @(8, 1.48, 317883.0),
    @(9, 0.82, 47584.0),
    @(10, 1.26, 65531.0),
    @(11, 0.87, 65396.0),
    @(12, 1.2, 190155.0),
    @(13, 0.45, 17980.0),
    @(14, 0.95, 67490.0),
    @(15, 1.23, 294098.0),
    @(16, 1.68, 198919.0),
    @(17, 0.96, 110707.0),
    @(19, 1.01, 106782.0),
    @(20, 1.33, 118407.0),
    @(21, 0.75, 69775.0),
    @(22, 1.74, 95494.0),
    @(23, 0.97, 55409.0),
    @(24, 2.33, 388564.0),
    @(25, 1.09, 218922.0),
    @(26, 1.62, 131507.0),
    @(27, 1.36, 146379.0),
    @(28, 0.51, 16781.0),
    @(29, 0.45, 19797.0),
    @(30, 0.67, 31823.0),
    @(31, 1.46, 106045.0),
    @(32, 0.97, 219155.0),
    @(33, 0.96, 97598.0),
    @(34, 1.55, 181511.0),
    @(35, 0.77, 46595.0),
    @(36, 2.17, 227953.0),
    @(37, 1.35, 71726.0),
    @(38, 1.07, 60061.0),
    @(39, 1.44, 62313.0),
    @(40, 1.8, 279389.0),
    @(41, 0.93, 58936.0),
    @(42, 1.17, 141360.0),
    @(43, 1.17, 61172.0),
    @(44, 3.17, 356919.0),
    @(45, 0.66, 37565.0),
    @(46, 1.55, 130940.0),
    @(47, 2.33, 288280.0)
)

foreach ($entry in $rowData) {
    $row = $entry[0]
    $eValue = $entry[1]
    $fValue = $entry[2]
    $ws.Range("E$row").Value = $eValue
    $ws.Range("F$row").Value = $fValue
}
